# Requirements Stack spreadsheet update
# - Bumps a few "Priority" (E column) values on the last rows of the table
# - Adjusts a handful of row heights (the wrapped-text rows grew a bit taller,
#   previously-default rows picked up an explicit height)
# - Leaves the view scrolled/zoomed in on the bottom of the table with F32 selected

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Cell value changes -------------------------------------------------
$ws.Range("E32").Value = 4

$ws.Range("E36").Value = 4
$ws.Range("E37").Value = 4
$ws.Range("E38").Value = 4
$ws.Range("E39").Value = 4
$ws.Range("E40").Value = 4

# ---- Row height changes --------------------------------------------------
# Rows that grow from the wrap-driven 29.25 height to 32.1
$tallerRows = 3,17,18,23,24,25,26,32,37,38,39,40
foreach ($r in $tallerRows) {
    $ws.Rows.Item($r).RowHeight = 32.1
}

# Rows that pick up (or keep) the shorter explicit 15.95 height
$shortRows = 2,4,5,6,9,10,11,12,13,19,20,21,22,27,28
foreach ($r in $shortRows) {
    $ws.Rows.Item($r).RowHeight = 15.95
}

# ---- View / window state -------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 2
$win.Zoom = 160

$ws.Range("F32").Select()
